# Update the "nombre_aides" (column C) and "montant_total" (column E) values
# for the regional rows affected by the 2022-06-15 Fonds de solidarite data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Row = 63; C = 40865; E = 116846014 }
    @{ Row = 111; C = 6011; E = 12081251 }
    @{ Row = 117; C = 19724; E = 56614531 }
    @{ Row = 125; C = 4603; E = 13149372 }
    @{ Row = 132; C = 6668; E = 13493261 }
    @{ Row = 134; C = 5681; E = 17175396 }
    @{ Row = 152; C = 126052; E = 716082530 }
    @{ Row = 168; C = 285071; E = 1212579919 }
    @{ Row = 169; C = 562644; E = 1285341892 }
    @{ Row = 170; C = 367484; E = 2847262741 }
    @{ Row = 171; C = 115200; E = 448106930 }
    @{ Row = 172; C = 21639; E = 73065656 }
    @{ Row = 174; C = 357315; E = 1019595846 }
    @{ Row = 175; C = 125587; E = 814144515 }
    @{ Row = 178; C = 75367; E = 102782703 }
    @{ Row = 179; C = 235759; E = 813180199 }
    @{ Row = 180; C = 141509; E = 341141722 }
    @{ Row = 249; C = 37151; E = 148715895 }
    @{ Row = 257; C = 182552; E = 1063829900 }
    @{ Row = 273; C = 88873; E = 139963440 }
    @{ Row = 295; C = 91338; E = 552983184 }
    @{ Row = 299; C = 43284; E = 140797061 }
    @{ Row = 311; C = 190857; E = 586450899 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

$wb.Save()
